# Regenerate s_vals data to filter save games.
# Updates columns B (TB), C (d2S), D (K), E (IP), G (sum) for rows 2-15.
# Column F (Win) and column A (dates) are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @{
    2  = @{ B = 0.1554434735375247;  C = 0.3375848360084654;  D = 3.082599426703578;  E = 0.4998867070740569;  G = 4.075514443323626 }
    3  = @{ B = 3.182878228561681;   C = 1.65323645889881;    D = 0.7127328510149897; E = 0.4998867070740569;  G = 6.048734245549538 }
    4  = @{ B = 3.182878228561681;   C = 1.65323645889881;    D = 0.7127328510149897; E = 0.4998867070740569;  G = 6.048734245549538 }
    5  = @{ B = 3.182878228561681;   C = 1.65323645889881;    D = 0.1529057820181812; E = 0.4998867070740569;  G = 5.488907176552729 }
    6  = @{ B = 0.7287194209349384;  C = 1.65323645889881;    D = 0.1529057820181812; E = 0.4998867070740569;  G = 3.034748368925986 }
    7  = @{ B = 3.182878228561681;   C = 9.226618575922256;   D = 157.8057217802531;  E = 6.48142807727062;    G = 176.6966466620077 }
    8  = @{ B = 3.182878228561681;   C = 1.65323645889881;    D = 0.1529057820181812; E = 0.4998867070740569;  G = 5.488907176552729 }
    9  = @{ B = 0.1554434735375247;  C = 9.226618575922256;   D = 0.1529057820181812; E = 6.48142807727062;    G = 16.01639590874858 }
    10 = @{ B = 3.182878228561681;   C = 1.65323645889881;    D = 0.7127328510149897; E = 0.4998867070740569;  G = 6.048734245549538 }
    11 = @{ B = 3.182878228561681;   C = 1.65323645889881;    D = 0.7127328510149897; E = 0.4998867070740569;  G = 6.048734245549538 }
    12 = @{ B = 0.7287194209349384;  C = 1.65323645889881;    D = 0.7127328510149897; E = 6.48142807727062;    G = 9.576116808119359 }
    13 = @{ B = 3.182878228561681;   C = 1.65323645889881;    D = 0.1529057820181812; E = 0.4998867070740569;  G = 5.488907176552729 }
    14 = @{ B = 3.182878228561681;   C = 1.65323645889881;    D = 0.7127328510149897; E = 0.4998867070740569;  G = 6.048734245549538 }
    15 = @{ B = 3.182878228561681;   C = 1.65323645889881;    D = 0.1529057820181812; E = 0.4998867070740569;  G = 5.488907176552729 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
